$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Leakage Modelling")

# Header G5: "Requirements of Data" (was pointing at a stale shared-string slot)
$ws.Range("G5").Value = "Requirements of Data"

# Row 6: existing reference entry - values re-aligned under the corrected headers
$ws.Range("C6").Value = "Optimization Tool to Improve the Management of the Leakages and Recovered Energy in Irrigation Water Systems"
$ws.Range("D6").Value = "Elsevier - Agricultural Water Management"
$ws.Range("E6").Value = "'2021"
$ws.Range("F6").Value = "Leakage modelling was based on the balance of water volumes and distributed along the network with the global emitter coefficient (K) concept."
$ws.Range("G6").Value = "Flow meter sensors."

# Row 7: brand new reference entry being added
$ws.Range("C7").Value = "Leakage Detection in Water Networks by a Calibration Method"
$ws.Range("D7").Value = "Elsevier - Flow Measurement and Instrumentation"
$ws.Range("E7").Value = "'2021"
$ws.Range("F7").Value = "The zone with the most leakage was identified by analysing the result of pressure sensors. Node demands and pipe roughness were calibrated with a metaheuristic optimization algorithm. Then the probability of leakage in each sub-zone was estimated."
$ws.Range("G7").Value = "Pressure and flow meter sensors."

# Re-apply the original "Requirements of Data" column formatting (quote-prefixed
# text style) that Excel's plain .Value assignment above would otherwise reset.
$ws.Range("H6").Copy()
$ws.Range("G6").PasteSpecial(-4122)
$ws.Range("G7").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 7 grows taller to fit the new wrapped text
$ws.Rows.Item(7).RowHeight = 100.8

# Reset the view: no frozen/offscreen top-left cell, new selection on row 8
$ws.Activate()
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("C8:F8").Select()
